$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing Difficulty cells from "Easy" to "Medium" ---
$ws.Range("E3").Value = "Medium"
$ws.Range("E5").Value = "Medium"
$ws.Range("F6").Value = "Medium"
$ws.Range("E7").Value = "Medium"
$ws.Range("E8").Value = "Medium"

# --- Add new row 9: "Merge Two Sorted Lists" ---
$ws.Range("A9").Value = "Merge Two Sorted Lists"
$ws.Range("B9").Value = "Linked List"
$ws.Range("C9").Value = "Yes"
$ws.Range("D9").Value = "No"
$ws.Range("E9").Value = "Easy"
$ws.Range("F9").Value = "Easy"
$ws.Hyperlinks.Add($ws.Range("G9"), "21%20-%20Merge%20Two%20Sorted%20Lists", "", "", "21 - Merge Two Sorted Lists")
$ws.Range("G9").Style = "Hyperlink"

# --- Extend conditional formatting to cover the new row ---
$cf = $ws.Range("D9:F9").FormatConditions
$cf.Add(1, 3, '"Hard"')
$cf.Add(1, 3, '"Medium"')
$cf.Add(1, 3, '"Easy"')

# --- Extend data validation ranges to include row 9 (recreated in original order) ---
$ws.Range("E2:F9").Validation.Delete()
$ws.Range("E2:F9").Validation.Add(3, 1, 1, '"Easy, Medium, Hard"')

$ws.Range("C2:C9").Validation.Delete()
$ws.Range("C2:C9").Validation.Add(3, 1, 1, '"Yes, No"')
$ws.Range("C2:C9").Validation.IgnoreBlank = $false

$ws.Range("B2:B9").Validation.Delete()
$ws.Range("B2:B9").Validation.Add(3, 1, 1, '"Array, Binary, Dynamic Programming, Graph, Interval, Linked List, Matrix, String, Tree, Heap"')

$ws.Range("D2:D9").Validation.Delete()
$ws.Range("D2:D9").Validation.Add(3, 1, 1, '"Yes, No"')

# --- Update selection to match the saved view state ---
$ws.Range("G16").Select() | Out-Null
